$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-16 Wednesday", "2025-07-17 Thursday"),
    @("236×9=", "872×3="),
    @("781×6=", "774×6="),
    @("470×9=", "694×3="),
    @("829×3=", "214×2="),
    @("925×8=", "359×8="),
    @("540×8=", "974×9="),
    @("824×6=", "294×8="),
    @("897×9=", "829×2="),
    @("106×5=", "387×5="),
    @("215×6=", "404×5="),
    @("253×4=", "337×7="),
    @("515×9=", "272×8="),
    @("290×8=", "507×8="),
    @("612×6=", "359×3="),
    @("870×7=", "289×2="),
    @("597×5=", "199×5="),
    @("450×8=", "132×4="),
    @("651×4=", "875×4="),
    @("611×7=", "312×5="),
    @("131×4=", "828×4="),
    @("974×7=", "125×6="),
    @("486×2=", "838×6="),
    @("884×8=", "635×6="),
    @("632×6=", "755×8="),
    @("841×4=", "656×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
